$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1, cell A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 08:08"

# Row 6 - India
$ws.Range("B6").Value = 698233
$ws.Range("C6").Value = 397
$ws.Range("D6").Value = 424928
$ws.Range("E6").Value = 253602
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 19703

# Row 47 - Afganistan
$ws.Range("B47").Value = 33190
$ws.Range("C47").Value = 239
$ws.Range("D47").Value = 20103
$ws.Range("E47").Value = 12189
$ws.Range("G47").Value = 34
$ws.Range("H47").Value = 898

# Row 71 - Uzbekistan
$ws.Range("B71").Value = 10143
$ws.Range("C71").Value = 123
$ws.Range("E71").Value = 3524
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 35

# Row 75 - Australia
$ws.Range("B75").Value = 8586
$ws.Range("C75").Value = 137
$ws.Range("E75").Value = 1060

# Row 90 - Costa Rica
$ws.Range("E90").Value = 3231
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 20

# Row 100 - Tailandia
$ws.Range("B100").Value = 3195
$ws.Range("C100").Value = 5
$ws.Range("D100").Value = 3072
$ws.Range("E100").Value = 65

# Row 141 - Georgia
$ws.Range("B141").Value = 953
$ws.Range("C141").Value = 2
$ws.Range("D141").Value = 830
